$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("D9").Value = "데이터 사이언스 스터디 교재 추천 + alpha"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/study-textbook-recommendation/#utm_source=rss&utm_medium=rss&utm_campaign=study-textbook-recommendation"

# Row 39
$ws.Range("D39").Value = "Deep Face Recognition with ArcFace in Keras and Python"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Deep-Face-Recognition-with-ArcFace-in-Keras-and-Python-1"

# Row 41
$ws.Range("D41").Value = "Service Mesh 에서의 Sidecar"
$ws.Range("E41").Value = "http://cloudinsight.net/cloud/service-mesh-%ec%97%90%ec%84%9c%ec%9d%98-sidecar/"

# Row 46
$ws.Range("D46").Value = "심장전도시스템 (cardiac conduction system) ②"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/376"
